# no-op test
$d = $word.ActiveDocument
Write-Host "Test"
